# Update the forward-flight data table and housekeeping (selection, cell
# styles) to match the latest data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Updated Thrust Required values (column A). Column B - speed (kt) -
# is unchanged. ---
$newValues = @{
    2  = 178.30180551368241
    3  = 167.9995536168438
    4  = 159.38189133788265
    5  = 152.20245354675555
    6  = 146.26132249585379
    7  = 141.39490264339935
    8  = 137.46828002633853
    9  = 134.36939308186311
    10 = 132.00454022650871
    11 = 130.29488499409052
    12 = 129.17371337918624
    13 = 128.58426389184271
    14 = 128.47799761590332
    15 = 128.81320918193796
    16 = 129.55390398115915
    17 = 130.66888485625759
    18 = 132.13100476613323
    19 = 133.91655182758191
    20 = 136.00474059870376
    21 = 138.37728913359408
    22 = 141.01806567055226
    23 = 143.91279215324514
    24 = 147.0487943719256
    25 = 150.41479053103544
    26 = 154.00071163471534
    27 = 157.79754833336509
    28 = 161.79721986806337
    29 = 165.99246154266839
    30 = 170.37672778945876
    31 = 174.9441084067615
    32 = 179.68925596199932
    33 = 184.60732269103482
    34 = 189.69390550024062
    35 = 194.94499790364128
    36 = 200.35694791341913
    37 = 205.92642105569158
    38 = 211.65036781083293
    39 = 217.52599488358027
    40 = 223.55073979661577
    41 = 229.72224837539554
    42 = 236.03835475422156
    43 = 242.49706358598283
    44 = 249.09653418229473
    45 = 255.83506634830687
    46 = 262.71108770834684
    47 = 269.72314234573821
    48 = 276.86988060333903
    49 = 284.15004991121532
    50 = 291.56248652491348
    51 = 299.1061080724611
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 1).Value = $newValues[$row]
}

# --- Rows whose direct formatting is cleared in the refreshed sheet (they
# fall back to plain/default formatting instead of the "Note" highlight
# they used to carry). ---
$rowsToClear = @(16, 19, 20, 21, 22, 23, 25, 26, 27, 28, 29, 30, 39, 40, 42)
foreach ($row in $rowsToClear) {
    $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 2)).ClearFormats()
}

# --- Rows that now carry the "Note" highlight style (row 41 previously used
# the "Good" style, which is being retired). ---
$rowsToNote = @(24, 41)
foreach ($row in $rowsToNote) {
    $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 2)).Style = "Note"
}

# The "Good" cell style is no longer used anywhere in the workbook - drop it.
$wb.Styles("Good").Delete()

# --- Selection moves to K15 (was C1:C1048576). ---
[void]$ws.Range("K15").Select()

# --- Window position moved on screen. ---
try {
    $win = $wb.Windows.Item(1)
    $win.Left = 7500
    $win.Top = 460
} catch {}
